# Add a new "Save" column (H) to the sheet, mirroring the header style
# used by the existing columns, and populate the per-row save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting (bold font, border, centered alignment)
# from the neighboring "sum" header G1, then overwrite its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Per-row "Save" values for rows 2-14. All zero except row 11, which is 1.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
